$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "dnasr281@gmail.com, *") {
        $parts = $val -split ", ", 2
        if ($parts.Count -eq 2) {
            $cell.Value = "$($parts[1]), $($parts[0])"
        }
    }
}
